$d = $word.ActiveDocument

$replacements = @(
    @("981×4=", "902×7="),
    @("809×8=", "608×8="),
    @("725×5=", "442×7="),
    @("598×3=", "509×8="),
    @("986×4=", "162×4="),
    @("275×7=", "776×8="),
    @("244×7=", "205×7="),
    @("948×8=", "194×2="),
    @("120×7=", "743×9="),
    @("453×4=", "993×8="),
    @("949×3=", "329×4="),
    @("883×4=", "438×6="),
    @("682×8=", "441×8="),
    @("137×3=", "853×5="),
    @("639×5=", "341×9="),
    @("793×8=", "370×7="),
    @("597×4=", "890×4="),
    @("828×5=", "678×4="),
    @("963×3=", "380×3="),
    @("738×4=", "447×3="),
    @("548×4=", "583×4="),
    @("649×9=", "317×7="),
    @("527×4=", "849×2="),
    @("860×8=", "887×9="),
    @("433×9=", "486×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
